$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove offensive/stray notes left in column H (rows 1 and 2), keep H3 ("cocksmal") as-is.
$ws.Range("H1").ClearContents() | Out-Null
$ws.Range("H2").ClearContents() | Out-Null

# Restore cursor position to where it was left after the edit.
$ws.Range("M12").Select() | Out-Null
$ws.Range("W14").Select() | Out-Null
